# Commit: "Updated with 1c task as completed And started with 2a"
#
# Task "1c" (row 10, column H - Status) is marked as completed.
# Task "2a" (row 13, column H - Status) is started: Start date / End date
# set to 29-Jul-2015 (serial 42214) and Status set to "in progress".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("iOS_Estimate")

# --- 1c: mark as completed ---
$ws.Range("H10").Value = "completed"

# --- 2a: started ---
$ws.Range("F13").Value = 42214
$ws.Range("F13").NumberFormat = "d-mmm"

$ws.Range("G13").Value = 42214
$ws.Range("G13").NumberFormat = "d-mmm"

$ws.Range("H13").Value = "in progress"

# Update the active cell / selection shown in the sheet view
$ws.Range("I13").Select()
